$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-08 Tuesday" "2024-10-09 Wednesday"

Replace-Text "439÷7=62, 5" "100÷9=11, 1"
Replace-Text "750÷4=187, 2" "222÷8=27, 6"
Replace-Text "882÷3=294, 0" "554÷4=138, 2"
Replace-Text "602÷8=75, 2" "147÷4=36, 3"
Replace-Text "139÷4=34, 3" "467÷5=93, 2"

Replace-Text "810÷5=162, 0" "519÷2=259, 1"
Replace-Text "369÷7=52, 5" "102÷8=12, 6"
Replace-Text "252÷6=42, 0" "450÷2=225, 0"
Replace-Text "222÷3=74, 0" "431÷6=71, 5"
Replace-Text "617÷7=88, 1" "774÷6=129, 0"

Replace-Text "838÷8=104, 6" "271÷4=67, 3"
Replace-Text "278÷5=55, 3" "918÷4=229, 2"
Replace-Text "581÷6=96, 5" "514÷9=57, 1"
Replace-Text "890÷2=445, 0" "285÷3=95, 0"
Replace-Text "959÷6=159, 5" "880÷8=110, 0"

Replace-Text "459÷3=153, 0" "756÷6=126, 0"
Replace-Text "716÷8=89, 4" "843÷3=281, 0"
Replace-Text "354÷8=44, 2" "610÷4=152, 2"
Replace-Text "831÷4=207, 3" "853÷7=121, 6"
Replace-Text "545÷9=60, 5" "802÷6=133, 4"

Replace-Text "278÷3=92, 2" "484÷8=60, 4"
Replace-Text "495÷6=82, 3" "392÷4=98, 0"
Replace-Text "363÷5=72, 3" "182÷4=45, 2"
Replace-Text "535÷7=76, 3" "752÷4=188, 0"
Replace-Text "600÷4=150, 0" "775÷8=96, 7"

Write-Output "done"
